$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated cluster-assignment fractions (rows 2-35)
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.008635578583765112
$ws.Range("F2").Value = 0.008866615265998464
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.02245250431778929
$ws.Range("F3").Value = 0.08840914931894027
$ws.Range("G3").Value = 0
$ws.Range("F4").Value = 0.01285016705217171
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.1206896551724138
$ws.Range("B5").Value = 0.3846153846153846
$ws.Range("C5").Value = 0.5440414507772044
$ws.Range("H5").Value = 0.15
$ws.Range("I5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0.01381692573402418
$ws.Range("F6").Value = 0.01297866872269343
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0.001028013364173734
$ws.Range("G7").Value = 0.05812291567413057
$ws.Range("J7").Value = 0.04545454545454546
$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 0.1714285714285714
$ws.Range("F8").Value = 0.08866615265998369
$ws.Range("G8").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0.001028013364173734
$ws.Range("G9").Value = 0.1357789423535017
$ws.Range("J9").Value = 0.3181818181818182
$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 0.06622201048118155
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0.2241379310344828
$ws.Range("E11").Value = 0
$ws.Range("G11").Value = 0.01238685088137209
$ws.Range("E12").Value = 0
$ws.Range("G12").Value = 0.02858504049547405
$ws.Range("I12").Value = 0.1206896551724138
$ws.Range("B13").Value = 0.02797202797202797
$ws.Range("C13").Value = 0.06390328151986178
$ws.Range("F13").Value = 0.002441531739912619
$ws.Range("B14").Value = 0.02097902097902098
$ws.Range("C14").Value = 0
$ws.Range("H14").Value = 0.8500000000000005
$ws.Range("I14").Value = 0
$ws.Range("B16").Value = 0.04895104895104895
$ws.Range("C16").Value = 0.07599309153713292
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0.06772038036494485
$ws.Range("G16").Value = 0.0738446879466413
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0.01724137931034483
$ws.Range("F17").Value = 0.01285016705217171
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 0.008099094807050979
$ws.Range("I18").Value = 0.01724137931034483
$ws.Range("B19").Value = 0.0979020979020979
$ws.Range("C19").Value = 0.09844559585492217
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0.03626943005181347
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0.05461320997172994
$ws.Range("G20").Value = 0.007146260123868511
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0.06044905008635574
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0.03919300950912382
$ws.Range("G21").Value = 0.08194378275369228
$ws.Range("B23").Value = 0.01398601398601399
$ws.Range("C23").Value = 0.04317789291882555
$ws.Range("F23").Value = 0.1072988948856317
$ws.Range("G23").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("G24").Value = 0.01381610290614579
$ws.Range("I24").Value = 0.08620689655172414
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0.002313030069390902
$ws.Range("G25").Value = 0.05097665555026206
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 0.0004764173415912339
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 0.02572653644592665
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0.005782575173477255
$ws.Range("G28").Value = 0.03239637922820392
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0.001670521716782318
$ws.Range("G29").Value = 0.02429728442115294
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0.002698535080956052
$ws.Range("G30").Value = 0.04144830871843738
$ws.Range("E31").Value = 0
$ws.Range("G31").Value = 0.05383515959980947
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0.01381692573402418
$ws.Range("F32").Value = 0.07465947057311723
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0.08620689655172414
$ws.Range("E33").Value = 0
$ws.Range("G33").Value = 0.1738923296808005
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0.3275862068965517
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 0.008635578583765112
$ws.Range("F34").Value = 0.004883063479825238
$ws.Range("E35").Value = 0
$ws.Range("G35").Value = 0.02953787517865652

# Remove the now-unused "Joint regime area" rows (36-40); sheet shrinks to A1:J35
$ws.Rows("36:40").Delete()
